# Update the "funded" grants sheet:
#  - Item 7 (Treating Complex Syntax in Children with DLD): update the USU
#    portion amount text on row 36.
#  - Add a brand-new item 8: the "Alzheimer's Disease in Native Hawaiians
#    and Pacific Islanders..." NIH U01 grant (Dr. Meek / Perry Ridge, BYU
#    subcontract) that was funded, filling rows 38-41.
#  - Make "funded" the active sheet/tab (it had been "dead").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("funded")

# --- Item 7: amend the amount note on row 36 ---
$ws.Range("F36").Value = "Amount\`$1,266,865 (USU's portion)"

# --- Item 8 (new row block, rows 38-41) ---
$ws.Range("A38").Value = 8
$ws.Range("B38").Value = "NIH: U01"
$ws.Range("C38").Value = "2025-2030"
$ws.Range("D38").Value = "Alzheimer's Disease in Native Hawaiians and Pacific Islanders: Sample Acquisition and Whole Genome Sequencing"
$ws.Range("E38").Value = "Subcontract to BYU"
$ws.Range("F38").Value = "PI: Perry Ridge (BYU-P)"

$ws.Range("F39").Value = "co-PIs: John Kauwe (BYU-H), JoAnn Tschanz (USU)"

$ws.Range("F40").Value = "Roll: statistician"

$ws.Range("F41").Value = "Amount \`$2,058,829 (USU's portion)"

# --- Make "funded" the active tab / selected cell ---
$ws.Activate()
[void]$ws.Range("F45").Select()
